# Update column F (dSF) values on rows 2-15 per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -2
    4  = 3
    5  = -4
    6  = -4
    7  = -4
    8  = 3
    9  = 9
    10 = -6
    11 = 6
    12 = -1
    13 = -2
    14 = -4
    15 = -4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
